$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.746.06"
$ws.Range("E2").Value = "  +5.75%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.432.27"
$ws.Range("E3").Value = "  +6.99%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.01"
$ws.Range("E5").Value = "  +7.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.38"
$ws.Range("E6").Value = "  +7.61%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.438.64"
$ws.Range("E8").Value = "  +6.93%  "

$ws.Range("E9").Value = "  +0.80%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.59"
$ws.Range("E10").Value = "  +3.21%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.123"
$ws.Range("E11").Value = "  +8.65%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.438"
$ws.Range("E12").Value = "  +1.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.024.65"
$ws.Range("E13").Value = "  +7.02%  "

$ws.Range("E14").Value = "  -0.67%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000187"
$ws.Range("E15").Value = "  +8.39%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.34"
$ws.Range("E16").Value = "  +5.58%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.821.91"
$ws.Range("E17").Value = "  +5.88%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.432.92"
$ws.Range("E18").Value = "  +7.14%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.44"
$ws.Range("E19").Value = "  +2.56%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.32"
$ws.Range("E20").Value = "  +7.78%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.50"
$ws.Range("E21").Value = "  +3.13%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "391.53"
$ws.Range("E22").Value = "  +5.19%  "

$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.541"
$ws.Range("E23").Value = "  +3.34%  "

$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.996"
$ws.Range("E24").Value = "  -0.53%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.08"
$ws.Range("E25").Value = "  +3.36%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000109"
$ws.Range("E26").Value = "  +23.48%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.53"
$ws.Range("E27").Value = "  +9.80%  "

$ws.Range("E28").Value = "  +6.49%  "

$ws.Range("E29").Value = "  -0.03%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.72"
$ws.Range("E30").Value = "  +9.54%  "

$ws.Range("E31").Value = "  +15.93%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.05"
$ws.Range("E32").Value = "  +7.60%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.80"
$ws.Range("E33").Value = "  +9.73%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.54"
$ws.Range("E34").Value = "  +4.78%  "

$ws.Range("E35").Value = "  -0.05%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.85"
$ws.Range("E36").Value = "  +4.07%  "

$ws.Range("E37").Value = "  +9.08%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "158.66"
$ws.Range("E38").Value = "  -0.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "28.17"
$ws.Range("E39").Value = "  +6.72%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0783"
$ws.Range("E40").Value = "  +10.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.88"
$ws.Range("E41").Value = "  +10.79%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.881.52"
$ws.Range("E42").Value = "  +3.10%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0321"
$ws.Range("E43").Value = "  +2.12%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.772"
$ws.Range("E44").Value = "  +7.22%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.93"
$ws.Range("E45").Value = "  +4.91%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.38"
$ws.Range("E46").Value = "  +3.67%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.10"
$ws.Range("E47").Value = "  +11.11%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.480.46"
$ws.Range("E48").Value = "  +7.11%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.65"
$ws.Range("E49").Value = "  +9.17%  "

$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.12"
$ws.Range("E50").Value = "  +25.61%  "

$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.40"
$ws.Range("E51").Value = "  +3.82%  "
